$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.238860845565796
$ws.Range("B1").Value = 2.451752901077271
$ws.Range("C1").Value = 4.79292106628418
$ws.Range("D1").Value = 3.122159719467163
$ws.Range("E1").Value = 1.133920907974243
